$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-25 Wednesday" "2026-02-26 Thursday"

Replace-Text "65÷7=9, 2" "84÷7=12, 0"
Replace-Text "70÷2=35, 0" "65÷5=13, 0"
Replace-Text "85÷5=17, 0" "54÷9=6, 0"
Replace-Text "49÷5=9, 4" "53÷7=7, 4"
Replace-Text "44÷5=8, 4" "53÷4=13, 1"
Replace-Text "87÷8=10, 7" "42÷7=6, 0"
Replace-Text "88÷3=29, 1" "80÷9=8, 8"
Replace-Text "14÷5=2, 4" "71÷7=10, 1"
Replace-Text "36÷2=18, 0" "65÷3=21, 2"
Replace-Text "70÷9=7, 7" "40÷7=5, 5"
Replace-Text "50÷6=8, 2" "33÷6=5, 3"
Replace-Text "98÷2=49, 0" "65÷5=13, 0"
Replace-Text "81÷3=27, 0" "89÷6=14, 5"
Replace-Text "14÷9=1, 5" "82÷6=13, 4"
Replace-Text "74÷6=12, 2" "54÷2=27, 0"
Replace-Text "92÷3=30, 2" "68÷5=13, 3"
Replace-Text "85÷9=9, 4" "11÷6=1, 5"
Replace-Text "53÷9=5, 8" "76÷7=10, 6"
Replace-Text "67÷9=7, 4" "90÷2=45, 0"
Replace-Text "81÷6=13, 3" "74÷2=37, 0"
Replace-Text "11÷3=3, 2" "90÷4=22, 2"
Replace-Text "21÷2=10, 1" "91÷7=13, 0"
Replace-Text "43÷6=7, 1" "28÷6=4, 4"
Replace-Text "16÷5=3, 1" "60÷4=15, 0"
Replace-Text "34÷6=5, 4" "48÷7=6, 6"
